$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextCell $ws 'D2' '29.057.33'
Set-TextCell $ws 'E2' '  +0.01%  '
Set-TextCell $ws 'D3' '1.835.76'
Set-TextCell $ws 'D4' '0.9989'
Set-TextCell $ws 'E4' '  +0.02%  '
Set-TextCell $ws 'D5' '244.17'
Set-TextCell $ws 'E5' '  +1.38%  '
Set-TextCell $ws 'D6' '0.6333'
Set-TextCell $ws 'E6' '  +2.07%  '
Set-TextCell $ws 'E7' '  +0.05%  '
Set-TextCell $ws 'D8' '0.07571'
Set-TextCell $ws 'E8' '  +2.93%  '
Set-TextCell $ws 'E9' '  +1.01%  '
Set-TextCell $ws 'D10' '22.82'
Set-TextCell $ws 'E10' '  +0.55%  '
Set-TextCell $ws 'D11' '0.07744'
Set-TextCell $ws 'E11' '  +0.79%  '
Set-TextCell $ws 'D12' '1.844.11'
Set-TextCell $ws 'E12' '  +0.82%  '
Set-TextCell $ws 'D13' '5.001'
Set-TextCell $ws 'E13' '  +0.73%  '
Set-TextCell $ws 'D14' '0.6713'
Set-TextCell $ws 'E14' '  +1.43%  '
Set-TextCell $ws 'D15' '83.25'
Set-TextCell $ws 'E15' '  +1.60%  '
Set-TextCell $ws 'D16' '0.000009806'
Set-TextCell $ws 'E16' '  +8.46%  '
Set-TextCell $ws 'D17' '6.124'
Set-TextCell $ws 'E17' '  +1.71%  '
Set-TextCell $ws 'D18' '29.088.22'
Set-TextCell $ws 'E18' '  +0.10%  '
Set-TextCell $ws 'D19' '12.57'
Set-TextCell $ws 'E19' '  +1.66%  '
Set-TextCell $ws 'D20' '227.10'
Set-TextCell $ws 'E20' '  +0.82%  '
Set-TextCell $ws 'D21' '1.0000'
Set-TextCell $ws 'E21' '  -0.03%  '
Set-TextCell $ws 'D22' '7.237'
Set-TextCell $ws 'E22' '  +1.30%  '
Set-TextCell $ws 'D23' '1.000'
Set-TextCell $ws 'E23' '  +0.01%  '
Set-TextCell $ws 'D24' '160.66'
Set-TextCell $ws 'E24' '  +0.58%  '
Set-TextCell $ws 'E25' '  +3.71%  '
Set-TextCell $ws 'D26' '8.552'
Set-TextCell $ws 'E26' '  +1.57%  '
Set-TextCell $ws 'D27' '17.99'
Set-TextCell $ws 'E27' '  +1.09%  '
Set-TextCell $ws 'D28' '1.501'
Set-TextCell $ws 'E28' '  +0.22%  '
Set-TextCell $ws 'D29' '4.122'
Set-TextCell $ws 'E29' '  +1.80%  '
Set-TextCell $ws 'D30' '4.058'
Set-TextCell $ws 'E30' '  +0.54%  '
Set-TextCell $ws 'D31' '1.202'
Set-TextCell $ws 'E31' '  +0.24%  '
Set-TextCell $ws 'D32' '0.05382'
Set-TextCell $ws 'E32' '  +2.70%  '
Set-TextCell $ws 'D33' '1.865'
Set-TextCell $ws 'D34' '0.7486'
Set-TextCell $ws 'E34' '  +2.39%  '
Set-TextCell $ws 'E35' '  -0.76%  '
Set-TextCell $ws 'D36' '2.675'
Set-TextCell $ws 'E36' '  +1.04%  '
Set-TextCell $ws 'D37' '1.247.15'
Set-TextCell $ws 'E38' '  +0.80%  '
Set-TextCell $ws 'D39' '2.760'
Set-TextCell $ws 'E39' '  +0.31%  '
Set-TextCell $ws 'D40' '6.609'
Set-TextCell $ws 'E40' '  +5.12%  '
Set-TextCell $ws 'D41' '0.9068'
Set-TextCell $ws 'E41' '  +0.66%  '
Set-TextCell $ws 'E42' '  +0.18%  '
Set-TextCell $ws 'D43' '102.79'
Set-TextCell $ws 'E43' '  +0.89%  '
Set-TextCell $ws 'D44' '1.989.18'
Set-TextCell $ws 'E44' '  +0.69%  '
Set-TextCell $ws 'E45' '  +3.61%  '
Set-TextCell $ws 'D46' '64.91'
Set-TextCell $ws 'E46' '  +1.45%  '
Set-TextCell $ws 'D47' '0.5112'
Set-TextCell $ws 'E47' '  -0.05%  '
Set-TextCell $ws 'D48' '0.4100'
Set-TextCell $ws 'E48' '  +3.41%  '
Set-TextCell $ws 'D49' '9.092'
Set-TextCell $ws 'E49' '  +3.47%  '
Set-TextCell $ws 'D50' '0.05793'
Set-TextCell $ws 'E50' '  +0.09%  '
Set-TextCell $ws 'D51' '6.777'
Set-TextCell $ws 'E51' '  +1.77%  '
